# Add a "grid frequency" table variant without an emergency state:
# update the two frequency readings that previously held a placeholder 0,
# and move the sheet's active selection down to the newly edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# B28 (time 14:30) now reports a frequency of 1 instead of 0
$ws.Range("B28").Value = 1

# B36 (time 16:30) now reports a frequency of 8 instead of 0
$ws.Range("B36").Value = 8

# Move the active selection to B37, mirroring where the author left off
$ws.Range("B37").Select()
